$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A=ECs D=ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Mrc1"
$ws.Cells.Item(2, 3).Value = "Ptprc"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.316305
$ws.Cells.Item(2, 8).Value = 0.948915
$ws.Cells.Item(2, 9).Value = 0.002355313614469985
$ws.Cells.Item(2, 10).Value = 0.002355313614469984
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1546356666666667
$ws.Cells.Item(2, 14).Value = 0.463907
$ws.Cells.Item(2, 15).Value = 0.0006541814359458435
$ws.Cells.Item(2, 16).Value = 0.0006541814359458435
$ws.Cells.Item(2, 17).Value = 0.048912034545
$ws.Cells.Item(2, 18).Value = 0.440208310905
$ws.Cells.Item(2, 19).Value = 0.00000154080244241677
$ws.Cells.Item(2, 20).Value = 0.000001540802442416769

# Row 3: A=ECs D=FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Mrc1"
$ws.Cells.Item(3, 3).Value = "Ptprc"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.316305
$ws.Cells.Item(3, 8).Value = 0.948915
$ws.Cells.Item(3, 9).Value = 0.002355313614469985
$ws.Cells.Item(3, 10).Value = 0.002355313614469984
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3253436666666666
$ws.Cells.Item(3, 14).Value = 0.976031
$ws.Cells.Item(3, 15).Value = 0.001376356384162467
$ws.Cells.Item(3, 16).Value = 0.001376356384162467
$ws.Cells.Item(3, 17).Value = 0.102907828485
$ws.Cells.Item(3, 18).Value = 0.9261704563649999
$ws.Cells.Item(3, 19).Value = 0.000003241750929980539
$ws.Cells.Item(3, 20).Value = 0.000003241750929980538

# Row 4: A=ECs D=M2
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Mrc1"
$ws.Cells.Item(4, 3).Value = "Ptprc"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.316305
$ws.Cells.Item(4, 8).Value = 0.948915
$ws.Cells.Item(4, 9).Value = 0.002355313614469985
$ws.Cells.Item(4, 10).Value = 0.002355313614469984
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 231.9004236666667
$ws.Cells.Item(4, 14).Value = 695.701271
$ws.Cells.Item(4, 15).Value = 0.9810476161216116
$ws.Cells.Item(4, 16).Value = 0.9810476161216116
$ws.Cells.Item(4, 17).Value = 73.35126350788501
$ws.Cells.Item(4, 18).Value = 660.161371570965
$ws.Cells.Item(4, 19).Value = 0.002310674806694555
$ws.Cells.Item(4, 20).Value = 0.002310674806694554

# Row 5: A=ECs D=sCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Mrc1"
$ws.Cells.Item(5, 3).Value = "Ptprc"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.316305
$ws.Cells.Item(5, 8).Value = 0.948915
$ws.Cells.Item(5, 9).Value = 0.002355313614469985
$ws.Cells.Item(5, 10).Value = 0.002355313614469984
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.999992666666667
$ws.Cells.Item(5, 14).Value = 11.999978
$ws.Cells.Item(5, 15).Value = 0.01692184605828007
$ws.Cells.Item(5, 16).Value = 0.01692184605828007
$ws.Cells.Item(5, 17).Value = 1.26521768043
$ws.Cells.Item(5, 18).Value = 11.38695912387
$ws.Cells.Item(5, 19).Value = 0.0000398562544030323
$ws.Cells.Item(5, 20).Value = 0.00003985625440303228

# Row 6: A=FAPs D=ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Mrc1"
$ws.Cells.Item(6, 3).Value = "Ptprc"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.2150266666666667
$ws.Cells.Item(6, 8).Value = 0.64508
$ws.Cells.Item(6, 9).Value = 0.001601161016974437
$ws.Cells.Item(6, 10).Value = 0.001601161016974436
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1546356666666667
$ws.Cells.Item(6, 14).Value = 0.463907
$ws.Cells.Item(6, 15).Value = 0.0006541814359458435
$ws.Cells.Item(6, 16).Value = 0.0006541814359458435
$ws.Cells.Item(6, 17).Value = 0.03325079195111112
$ws.Cells.Item(6, 18).Value = 0.29925712756
$ws.Cells.Item(6, 19).Value = 0.000001047449813264844
$ws.Cells.Item(6, 20).Value = 0.000001047449813264844

# Row 7: A=FAPs D=FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Mrc1"
$ws.Cells.Item(7, 3).Value = "Ptprc"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.2150266666666667
$ws.Cells.Item(7, 8).Value = 0.64508
$ws.Cells.Item(7, 9).Value = 0.001601161016974437
$ws.Cells.Item(7, 10).Value = 0.001601161016974436
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3253436666666666
$ws.Cells.Item(7, 14).Value = 0.976031
$ws.Cells.Item(7, 15).Value = 0.001376356384162467
$ws.Cells.Item(7, 16).Value = 0.001376356384162467
$ws.Cells.Item(7, 17).Value = 0.06995756416444444
$ws.Cells.Item(7, 18).Value = 0.62961807748
$ws.Cells.Item(7, 19).Value = 0.000002203768187784834
$ws.Cells.Item(7, 20).Value = 0.000002203768187784834

# Row 8: A=FAPs D=M2
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Mrc1"
$ws.Cells.Item(8, 3).Value = "Ptprc"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.2150266666666667
$ws.Cells.Item(8, 8).Value = 0.64508
$ws.Cells.Item(8, 9).Value = 0.001601161016974437
$ws.Cells.Item(8, 10).Value = 0.001601161016974436
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 231.9004236666667
$ws.Cells.Item(8, 14).Value = 695.701271
$ws.Cells.Item(8, 15).Value = 0.9810476161216116
$ws.Cells.Item(8, 16).Value = 0.9810476161216116
$ws.Cells.Item(8, 17).Value = 49.86477509963112
$ws.Cells.Item(8, 18).Value = 448.78297589668
$ws.Cells.Item(8, 19).Value = 0.001570815198729627
$ws.Cells.Item(8, 20).Value = 0.001570815198729626

# Row 9: A=FAPs D=sCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Mrc1"
$ws.Cells.Item(9, 3).Value = "Ptprc"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.2150266666666667
$ws.Cells.Item(9, 8).Value = 0.64508
$ws.Cells.Item(9, 9).Value = 0.001601161016974437
$ws.Cells.Item(9, 10).Value = 0.001601161016974436
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.999992666666667
$ws.Cells.Item(9, 14).Value = 11.999978
$ws.Cells.Item(9, 15).Value = 0.01692184605828007
$ws.Cells.Item(9, 16).Value = 0.01692184605828007
$ws.Cells.Item(9, 17).Value = 0.8601050898044446
$ws.Cells.Item(9, 18).Value = 7.74094580824
$ws.Cells.Item(9, 19).Value = 0.00002709460024376058
$ws.Cells.Item(9, 20).Value = 0.00002709460024376057

# Row 10: A=M2 D=ECs
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Mrc1"
$ws.Cells.Item(10, 3).Value = "Ptprc"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 130.539482
$ws.Cells.Item(10, 8).Value = 391.618446
$ws.Cells.Item(10, 9).Value = 0.9720409705204137
$ws.Cells.Item(10, 10).Value = 0.9720409705204136
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1546356666666667
$ws.Cells.Item(10, 14).Value = 0.463907
$ws.Cells.Item(10, 15).Value = 0.0006541814359458435
$ws.Cells.Item(10, 16).Value = 0.0006541814359458435
$ws.Cells.Item(10, 17).Value = 20.18605982539133
$ws.Cells.Item(10, 18).Value = 181.674538428522
$ws.Cells.Item(10, 19).Value = 0.0006358911578932356
$ws.Cells.Item(10, 20).Value = 0.0006358911578932356

# Row 11: A=M2 D=FAPs
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Mrc1"
$ws.Cells.Item(11, 3).Value = "Ptprc"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 130.539482
$ws.Cells.Item(11, 8).Value = 391.618446
$ws.Cells.Item(11, 9).Value = 0.9720409705204137
$ws.Cells.Item(11, 10).Value = 0.9720409705204136
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.3253436666666666
$ws.Cells.Item(11, 14).Value = 0.976031
$ws.Cells.Item(11, 15).Value = 0.001376356384162467
$ws.Cells.Item(11, 16).Value = 0.001376356384162467
$ws.Cells.Item(11, 17).Value = 42.47019371864733
$ws.Cells.Item(11, 18).Value = 382.231743467826
$ws.Cells.Item(11, 19).Value = 0.001337874795443252
$ws.Cells.Item(11, 20).Value = 0.001337874795443252

# Row 12: A=M2 D=M2
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Mrc1"
$ws.Cells.Item(12, 3).Value = "Ptprc"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 130.539482
$ws.Cells.Item(12, 8).Value = 391.618446
$ws.Cells.Item(12, 9).Value = 0.9720409705204137
$ws.Cells.Item(12, 10).Value = 0.9720409705204136
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 231.9004236666667
$ws.Cells.Item(12, 14).Value = 695.701271
$ws.Cells.Item(12, 15).Value = 0.9810476161216116
$ws.Cells.Item(12, 16).Value = 0.9810476161216116
$ws.Cells.Item(12, 17).Value = 30272.16118102721
$ws.Cells.Item(12, 18).Value = 272449.4506292448
$ws.Cells.Item(12, 19).Value = 0.9536184769015896
$ws.Cells.Item(12, 20).Value = 0.9536184769015895

# Row 13: A=M2 D=sCs
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Mrc1"
$ws.Cells.Item(13, 3).Value = "Ptprc"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 130.539482
$ws.Cells.Item(13, 8).Value = 391.618446
$ws.Cells.Item(13, 9).Value = 0.9720409705204137
$ws.Cells.Item(13, 10).Value = 0.9720409705204136
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.999992666666667
$ws.Cells.Item(13, 14).Value = 11.999978
$ws.Cells.Item(13, 15).Value = 0.01692184605828007
$ws.Cells.Item(13, 16).Value = 0.01692184605828007
$ws.Cells.Item(13, 17).Value = 522.1569707104653
$ws.Cells.Item(13, 18).Value = 4699.412736394188
$ws.Cells.Item(13, 19).Value = 0.0164487276654876
$ws.Cells.Item(13, 20).Value = 0.0164487276654876

# Row 14: A=sCs D=ECs
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Mrc1"
$ws.Cells.Item(14, 3).Value = "Ptprc"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 3.223404333333333
$ws.Cells.Item(14, 8).Value = 9.670213
$ws.Cells.Item(14, 9).Value = 0.02400255484814197
$ws.Cells.Item(14, 10).Value = 0.02400255484814196
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.1546356666666667
$ws.Cells.Item(14, 14).Value = 0.463907
$ws.Cells.Item(14, 15).Value = 0.0006541814359458435
$ws.Cells.Item(14, 16).Value = 0.0006541814359458435
$ws.Cells.Item(14, 17).Value = 0.4984532780212222
$ws.Cells.Item(14, 18).Value = 4.486079502191
$ws.Cells.Item(14, 19).Value = 0.00001570202579692638
$ws.Cells.Item(14, 20).Value = 0.00001570202579692638

# Row 15: A=sCs D=FAPs
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Mrc1"
$ws.Cells.Item(15, 3).Value = "Ptprc"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 3.223404333333333
$ws.Cells.Item(15, 8).Value = 9.670213
$ws.Cells.Item(15, 9).Value = 0.02400255484814197
$ws.Cells.Item(15, 10).Value = 0.02400255484814196
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.3253436666666666
$ws.Cells.Item(15, 14).Value = 0.976031
$ws.Cells.Item(15, 15).Value = 0.001376356384162467
$ws.Cells.Item(15, 16).Value = 0.001376356384162467
$ws.Cells.Item(15, 17).Value = 1.048714184955889
$ws.Cells.Item(15, 18).Value = 9.438427664603001
$ws.Cells.Item(15, 19).Value = 0.00003303606960144997
$ws.Cells.Item(15, 20).Value = 0.00003303606960144996

# Row 16: A=sCs D=M2
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Mrc1"
$ws.Cells.Item(16, 3).Value = "Ptprc"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 3.223404333333333
$ws.Cells.Item(16, 8).Value = 9.670213
$ws.Cells.Item(16, 9).Value = 0.02400255484814197
$ws.Cells.Item(16, 10).Value = 0.02400255484814196
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 231.9004236666667
$ws.Cells.Item(16, 14).Value = 695.701271
$ws.Cells.Item(16, 15).Value = 0.9810476161216116
$ws.Cells.Item(16, 16).Value = 0.9810476161216116
$ws.Cells.Item(16, 17).Value = 747.5088305489693
$ws.Cells.Item(16, 18).Value = 6727.579474940723
$ws.Cells.Item(16, 19).Value = 0.02354764921459791
$ws.Cells.Item(16, 20).Value = 0.0235476492145979

# Row 17: A=sCs D=sCs
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Mrc1"
$ws.Cells.Item(17, 3).Value = "Ptprc"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 3.223404333333333
$ws.Cells.Item(17, 8).Value = 9.670213
$ws.Cells.Item(17, 9).Value = 0.02400255484814197
$ws.Cells.Item(17, 10).Value = 0.02400255484814196
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3.999992666666667
$ws.Cells.Item(17, 14).Value = 11.999978
$ws.Cells.Item(17, 15).Value = 0.01692184605828007
$ws.Cells.Item(17, 16).Value = 0.01692184605828007
$ws.Cells.Item(17, 17).Value = 12.89359369503489
$ws.Cells.Item(17, 18).Value = 116.042343255314
$ws.Cells.Item(17, 19).Value = 0.0004061675381456823
$ws.Cells.Item(17, 20).Value = 0.0004061675381456822
